# Scheduled market-data refresh: update cached Leve profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns H:N) across the
# ALC / ARM / BSM / CRP / GSM / LTW / WVR sheets, per the latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3429
$ws.Range("I74").Value = 3200.75
$ws.Range("J74").Value = 3733.3333
$ws.Range("K74").Value = 3200.75
$ws.Range("L74").Value = 3733.3333
$ws.Range("M74").Value = -2264.75
$ws.Range("N74").Value = -5605.3333
$ws.Range("H77").Value = 3429
$ws.Range("I77").Value = 3200.75
$ws.Range("J77").Value = 3733.3333
$ws.Range("K77").Value = 16003.75
$ws.Range("L77").Value = 18666.6665
$ws.Range("M77").Value = -11323.75
$ws.Range("N77").Value = -28026.6665
$ws.Range("H121").Value = 456.42856
$ws.Range("J121").Value = 456.42856
$ws.Range("L121").Value = 1369.28568
$ws.Range("N121").Value = -4863.28568
$ws.Range("H125").Value = 10192998
$ws.Range("I125").Value = 610.6667
$ws.Range("J125").Value = 14015144
$ws.Range("K125").Value = 5496.0003
$ws.Range("L125").Value = 126136296
$ws.Range("M125").Value = -3036.0003
$ws.Range("N125").Value = -126141216
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 20409110
$ws.Range("I137").Value = 25000694
$ws.Range("J137").Value = 2074.889
$ws.Range("K137").Value = 75002082
$ws.Range("L137").Value = 6224.667
$ws.Range("M137").Value = -74999532
$ws.Range("N137").Value = -11324.667
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20507.865
$ws.Range("I32").Value = 3416.7637
$ws.Range("J32").Value = 255510.5
$ws.Range("K32").Value = 3416.7637
$ws.Range("L32").Value = 255510.5
$ws.Range("M32").Value = -3129.7637
$ws.Range("N32").Value = -256084.5
$ws.Range("H74").Value = 4517.8604
$ws.Range("I74").Value = 1419.8182
$ws.Range("J74").Value = 14741.4
$ws.Range("K74").Value = 1419.8182
$ws.Range("L74").Value = 14741.4
$ws.Range("M74").Value = -545.8181999999999
$ws.Range("N74").Value = -16489.4
$ws.Range("H77").Value = 4517.8604
$ws.Range("I77").Value = 1419.8182
$ws.Range("J77").Value = 14741.4
$ws.Range("K77").Value = 7099.090999999999
$ws.Range("L77").Value = 73707
$ws.Range("M77").Value = -2731.090999999999
$ws.Range("N77").Value = -82443
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7774.5884
$ws.Range("I86").Value = 2453.8572
$ws.Range("J86").Value = 11499.1
$ws.Range("K86").Value = 2453.8572
$ws.Range("L86").Value = 11499.1
$ws.Range("M86").Value = -1330.8572
$ws.Range("N86").Value = -13745.1
$ws.Range("H89").Value = 7774.5884
$ws.Range("I89").Value = 2453.8572
$ws.Range("J89").Value = 11499.1
$ws.Range("K89").Value = 12269.286
$ws.Range("L89").Value = 57495.5
$ws.Range("M89").Value = -6653.286
$ws.Range("N89").Value = -68727.5
$ws.Range("H134").Value = 22224910
$ws.Range("I134").Value = 27779790
$ws.Range("J134").Value = 5389.1113
$ws.Range("K134").Value = 83339370
$ws.Range("L134").Value = 16167.3339
$ws.Range("M134").Value = -83336835
$ws.Range("N134").Value = -21237.3339
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 23635.5
$ws.Range("J16").Value = 1170.4286
$ws.Range("L16").Value = 1170.4286
$ws.Range("N16").Value = -1744.4286
$ws.Range("H31").Value = 2225
$ws.Range("I31").Value = 1304.52
$ws.Range("K31").Value = 1304.52
$ws.Range("M31").Value = -1009.52
$ws.Range("H34").Value = 2225
$ws.Range("I34").Value = 1304.52
$ws.Range("K34").Value = 1304.52
$ws.Range("M34").Value = -1102.52
$ws.Range("H86").Value = 25004766
$ws.Range("I86").Value = 38467630
$ws.Range("J86").Value = 2301.1428
$ws.Range("K86").Value = 38467630
$ws.Range("L86").Value = 2301.1428
$ws.Range("M86").Value = -38466507
$ws.Range("N86").Value = -4547.1428
$ws.Range("H89").Value = 25004766
$ws.Range("I89").Value = 38467630
$ws.Range("J89").Value = 2301.1428
$ws.Range("K89").Value = 192338150
$ws.Range("L89").Value = 11505.714
$ws.Range("M89").Value = -192332534
$ws.Range("N89").Value = -22737.714
$ws.Range("H113").Value = 23635.5
$ws.Range("J113").Value = 1170.4286
$ws.Range("L113").Value = 1170.4286
$ws.Range("N113").Value = -5510.4286
$ws.Range("H132").Value = 2678.027
$ws.Range("I132").Value = 2092.3103
$ws.Range("J132").Value = 4801.25
$ws.Range("K132").Value = 6276.9309
$ws.Range("L132").Value = 14403.75
$ws.Range("M132").Value = -3746.9309
$ws.Range("N132").Value = -19463.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6097.871
$ws.Range("I70").Value = 6072
$ws.Range("K70").Value = 6072
$ws.Range("M70").Value = -5802
$ws.Range("H73").Value = 6097.871
$ws.Range("I73").Value = 6072
$ws.Range("K73").Value = 6072
$ws.Range("M73").Value = -5136
$ws.Range("H102").Value = 2222.4167
$ws.Range("I102").Value = 2495.5715
$ws.Range("J102").Value = 1840
$ws.Range("K102").Value = 2495.5715
$ws.Range("L102").Value = 1840
$ws.Range("M102").Value = -873.5715
$ws.Range("N102").Value = -5084
$ws.Range("H108").Value = 32500
$ws.Range("J108").Value = 32500
$ws.Range("L108").Value = 32500
$ws.Range("N108").Value = -40180
$ws.Range("H122").Value = 655114.0600000001
$ws.Range("I122").Value = 742162.0600000001
$ws.Range("J122").Value = 2254
$ws.Range("K122").Value = 2226486.18
$ws.Range("L122").Value = 6762
$ws.Range("M122").Value = -2224036.18
$ws.Range("N122").Value = -11662
$ws.Range("H132").Value = 3005.0417
$ws.Range("I132").Value = 2656.5134
$ws.Range("J132").Value = 4177.364
$ws.Range("K132").Value = 7969.540199999999
$ws.Range("L132").Value = 12532.092
$ws.Range("M132").Value = -5439.540199999999
$ws.Range("N132").Value = -17592.092
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2894.389
$ws.Range("I100").Value = 2442.7144
$ws.Range("J100").Value = 3181.818
$ws.Range("K100").Value = 2442.7144
$ws.Range("L100").Value = 3181.818
$ws.Range("M100").Value = -1901.7144
$ws.Range("N100").Value = -4263.818
$ws.Range("H132").Value = 7338.6
$ws.Range("I132").Value = 9753.154
$ws.Range("J132").Value = 5492.1763
$ws.Range("K132").Value = 29259.462
$ws.Range("L132").Value = 16476.5289
$ws.Range("M132").Value = -26729.462
$ws.Range("N132").Value = -21536.5289
$ws.Range("H136").Value = 2804.6462
$ws.Range("I136").Value = 1433.7736
$ws.Range("J136").Value = 8859.333000000001
$ws.Range("K136").Value = 4301.3208
$ws.Range("L136").Value = 26577.999
$ws.Range("M136").Value = -1751.3208
$ws.Range("N136").Value = -31677.999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 168014
$ws.Range("I122").Value = 251151
$ws.Range("J122").Value = 1740
$ws.Range("K122").Value = 753453
$ws.Range("L122").Value = 5220
$ws.Range("M122").Value = -751003
$ws.Range("N122").Value = -10120
$ws.Range("H132").Value = 9093291
$ws.Range("I132").Value = 14708165
$ws.Range("J132").Value = 2543.7144
$ws.Range("K132").Value = 44124495
$ws.Range("L132").Value = 7631.1432
$ws.Range("M132").Value = -44121965
$ws.Range("N132").Value = -12691.1432
$ws.Range("H136").Value = 19959.906
$ws.Range("I136").Value = 23793.953
$ws.Range("J136").Value = 3473.5
$ws.Range("K136").Value = 71381.859
$ws.Range("L136").Value = 10420.5
$ws.Range("M136").Value = -68831.859
$ws.Range("N136").Value = -15520.5
